$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 614.75
$ws.Range("I6").Value = 416.85715
$ws.Range("K6").Value = 1250.57145
$ws.Range("M6").Value = -1138.57145

# Row 39
$ws.Range("H39").Value = 369.6154
$ws.Range("I39").Value = 300.8889
$ws.Range("J39").Value = 524.25
$ws.Range("K39").Value = 902.6667
$ws.Range("L39").Value = 1572.75
$ws.Range("M39").Value = -606.6667
$ws.Range("N39").Value = -2164.75

# Row 45
$ws.Range("H45").Value = 3137.2856
$ws.Range("I45").Value = 769
$ws.Range("J45").Value = 4913.5
$ws.Range("K45").Value = 2307
$ws.Range("L45").Value = 14740.5
$ws.Range("M45").Value = -2115
$ws.Range("N45").Value = -15124.5

# Row 61
$ws.Range("H61").Value = 9596.666999999999
$ws.Range("I61").Value = 9596.666999999999
$ws.Range("K61").Value = 28790.001
$ws.Range("M61").Value = -28618.001

# Row 76
$ws.Range("H76").Value = 5382.619
$ws.Range("I76").Value = 4002.2856
$ws.Range("K76").Value = 4002.2856
$ws.Range("M76").Value = -3687.2856

# Row 79
$ws.Range("H79").Value = 5382.619
$ws.Range("I79").Value = 4002.2856
$ws.Range("K79").Value = 4002.2856
$ws.Range("M79").Value = -2910.2856

# Row 98
$ws.Range("H98").Value = 1327.4166
$ws.Range("I98").Value = 1261.4286
$ws.Range("K98").Value = 1261.4286
$ws.Range("M98").Value = 236.5714

# Row 122
$ws.Range("H122").Value = 1327.4166
$ws.Range("I122").Value = 1261.4286
$ws.Range("K122").Value = 3784.2858
$ws.Range("M122").Value = -1334.2858

# Row 129
$ws.Range("H129").Value = 2050.0435
$ws.Range("I129").Value = 939.36365
$ws.Range("J129").Value = 3068.1667
$ws.Range("K129").Value = 2818.09095
$ws.Range("L129").Value = 9204.500100000001
$ws.Range("M129").Value = 2181.90905
$ws.Range("N129").Value = -19204.5001

# Row 132
$ws.Range("H132").Value = 1518.5892
$ws.Range("I132").Value = 1394.4783
$ws.Range("K132").Value = 4183.4349
$ws.Range("M132").Value = -1653.4349

# Row 137
$ws.Range("H137").Value = 20836086
$ws.Range("I137").Value = 41668736
$ws.Range("J137").Value = 3436.375
$ws.Range("K137").Value = 125006208
$ws.Range("L137").Value = 10309.125
$ws.Range("M137").Value = -125003658
$ws.Range("N137").Value = -15409.125


$ws = $wb.Worksheets.Item("ARM")
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# Row 63
$ws.Range("H63").Value = 2568.7144
$ws.Range("I63").Value = 1996.8334
$ws.Range("K63").Value = 1996.8334
$ws.Range("M63").Value = -1310.8334

# Row 66
$ws.Range("H66").Value = 2568.7144
$ws.Range("I66").Value = 1996.8334
$ws.Range("K66").Value = 9984.166999999999
$ws.Range("M66").Value = -6552.166999999999

# Row 122
$ws.Range("H122").Value = 47621120
$ws.Range("I122").Value = 1411.5714
$ws.Range("J122").Value = 142860540
$ws.Range("K122").Value = 4234.7142
$ws.Range("L122").Value = 428581620
$ws.Range("M122").Value = -1784.7142
$ws.Range("N122").Value = -428586520

# Row 132
$ws.Range("H132").Value = 1999.3235
$ws.Range("I132").Value = 1503.4286
$ws.Range("K132").Value = 4510.2858
$ws.Range("M132").Value = -1980.2858

# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()


$ws = $wb.Worksheets.Item("BSM")
# Row 110
$ws.Range("H110").Value = 64998.332
$ws.Range("J110").Value = 64998.332
$ws.Range("L110").Value = 64998.332
$ws.Range("N110").Value = -73178.33199999999

# Row 134
$ws.Range("H134").Value = 3068.4062
$ws.Range("I134").Value = 1177
$ws.Range("J134").Value = 9823.429
$ws.Range("K134").Value = 3531
$ws.Range("L134").Value = 29470.287
$ws.Range("M134").Value = -996
$ws.Range("N134").Value = -34540.287


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 102714.82
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 102714.82
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 102714.82
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -103304.82

# Row 34
$ws.Range("H34").Value = 102714.82
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 102714.82
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 102714.82
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -103118.82

# Row 105
$ws.Range("H105").Value = 18205.5
$ws.Range("I105").Value = 17400
$ws.Range("K105").Value = 17400
$ws.Range("M105").Value = -15653

# Row 134
$ws.Range("H134").Value = 3943.4546
$ws.Range("I134").Value = 2608.925
$ws.Range("K134").Value = 7826.775000000001
$ws.Range("M134").Value = -5291.775000000001


$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 31379.938
$ws.Range("J7").Value = 100108.2
$ws.Range("L7").Value = 300324.6
$ws.Range("N7").Value = -300548.6

# Row 24
$ws.Range("H24").Value = 5
$ws.Range("I24").Value = 5
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 15
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 215
$ws.Range("N24").ClearContents()

# Row 113
$ws.Range("H113").Value = 83334776
$ws.Range("J113").Value = 100001336
$ws.Range("L113").Value = 300004008
$ws.Range("N113").Value = -300008348

# Row 117
$ws.Range("H117").Value = 1778.1177
$ws.Range("I117").Value = 1194.7142
$ws.Range("J117").Value = 2186.5
$ws.Range("K117").Value = 3584.1426
$ws.Range("L117").Value = 6559.5
$ws.Range("M117").Value = -142.1425999999997
$ws.Range("N117").Value = -13443.5

# Row 122
$ws.Range("H122").Value = 6482259.5
$ws.Range("I122").Value = 8547341
$ws.Range("J122").Value = 6212901
$ws.Range("K122").Value = 76926069
$ws.Range("L122").Value = 55916109
$ws.Range("M122").Value = -76923619
$ws.Range("N122").Value = -55921009


$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 4634.4165
$ws.Range("I126").Value = 1549.6666
$ws.Range("K126").Value = 4648.9998
$ws.Range("M126").Value = -2178.9998

# Row 132
$ws.Range("H132").Value = 3287.5417
$ws.Range("I132").Value = 2876.7834
$ws.Range("J132").Value = 5341.3335
$ws.Range("K132").Value = 8630.350199999999
$ws.Range("L132").Value = 16024.0005
$ws.Range("M132").Value = -6100.350199999999
$ws.Range("N132").Value = -21084.0005

# Row 140
$ws.Range("H140").Value = 66966.5
$ws.Range("J140").Value = 66966.5
$ws.Range("L140").Value = 66966.5
$ws.Range("N140").Value = -77326.5


$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5417
$ws.Range("I40").Value = 3583.3333
$ws.Range("J40").Value = 8167.5
$ws.Range("K40").Value = 3583.3333
$ws.Range("L40").Value = 8167.5
$ws.Range("M40").Value = -3447.3333
$ws.Range("N40").Value = -8439.5

# Row 61
$ws.Range("H61").Value = 4694.4736
$ws.Range("I61").Value = 2605.375
$ws.Range("K61").Value = 2605.375
$ws.Range("M61").Value = -2403.375

# Row 113
$ws.Range("H113").Value = 4694.4736
$ws.Range("I113").Value = 2605.375
$ws.Range("K113").Value = 2605.375
$ws.Range("M113").Value = -435.375

# Row 136
$ws.Range("H136").Value = 4940.7915
$ws.Range("I136").Value = 3095.4
$ws.Range("J136").Value = 14167.75
$ws.Range("K136").Value = 9286.200000000001
$ws.Range("L136").Value = 42503.25
$ws.Range("M136").Value = -6736.200000000001
$ws.Range("N136").Value = -47603.25

# Row 137
$ws.Range("H137").Value = 32000
$ws.Range("I137").Value = 32000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 32000
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -26900
$ws.Range("N137").ClearContents()

# Row 139
$ws.Range("H139").Value = 54253
$ws.Range("J139").Value = 69807.5
$ws.Range("L139").Value = 69807.5
$ws.Range("N139").Value = -80087.5


$ws = $wb.Worksheets.Item("WVR")
# Row 55
$ws.Range("H55").Value = 29992
$ws.Range("I55").Value = 22633
$ws.Range("J55").Value = 37351
$ws.Range("K55").Value = 22633
$ws.Range("L55").Value = 37351
$ws.Range("M55").Value = -22356
$ws.Range("N55").Value = -37905

# Row 62
$ws.Range("H62").Value = 4712.6
$ws.Range("I62").Value = 4526
$ws.Range("K62").Value = 4526
$ws.Range("M62").Value = -3902

# Row 65
$ws.Range("H65").Value = 4712.6
$ws.Range("I65").Value = 4526
$ws.Range("K65").Value = 22630
$ws.Range("M65").Value = -19510

# Row 93
$ws.Range("H93").Value = 55000
$ws.Range("J93").Value = 55000
$ws.Range("L93").Value = 55000
$ws.Range("N93").Value = -59992

# Row 132
$ws.Range("H132").Value = 4439.1753
$ws.Range("I132").Value = 2194.6743
$ws.Range("K132").Value = 6584.0229
$ws.Range("M132").Value = -4054.0229

# Row 136
$ws.Range("H136").Value = 2017.1482
$ws.Range("I136").Value = 1363.7693
$ws.Range("K136").Value = 4091.3079
$ws.Range("M136").Value = -1541.3079


